# Refresh the 2023 (column J) violent-crime counts to add the daily data
# for 2023-03-27 across the "Citywide Totals" summary sheet, the
# "By Neighborhood" summary sheet, and every affected per-neighborhood
# detail sheet. Some rows also saw a one-off correction to the 2022
# (column I) figure.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 1503
$ws.Range("I3").Value = 7488
$ws.Range("J3").Value = 1580
$ws.Range("I4").Value = 1757
$ws.Range("J4").Value = 358
$ws.Range("J5").Value = 110
$ws.Range("J6").Value = 2052
$ws.Range("I7").Value = 26201
$ws.Range("J7").Value = 5603

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 28
$ws.Range("J7").Value = 158
$ws.Range("J8").Value = 346
$ws.Range("J10").Value = 34
$ws.Range("J11").Value = 69
$ws.Range("J14").Value = 20
$ws.Range("J15").Value = 71
$ws.Range("J16").Value = 18
$ws.Range("J17").Value = 15
$ws.Range("J19").Value = 198
$ws.Range("J25").Value = 30
$ws.Range("J29").Value = 312
$ws.Range("J32").Value = 11
$ws.Range("J33").Value = 233
$ws.Range("J34").Value = 35
$ws.Range("J36").Value = 85
$ws.Range("J37").Value = 193
$ws.Range("J42").Value = 221
$ws.Range("J44").Value = 45
$ws.Range("J47").Value = 48
$ws.Range("J48").Value = 43
$ws.Range("J54").Value = 109
$ws.Range("J55").Value = 68
$ws.Range("J60").Value = 31
$ws.Range("I63").Value = 194
$ws.Range("J63").Value = 22
$ws.Range("J64").Value = 37
$ws.Range("J66").Value = 13
$ws.Range("J67").Value = 197
$ws.Range("J73").Value = 54
$ws.Range("J78").Value = 72
$ws.Range("J83").Value = 139
$ws.Range("J84").Value = 58
$ws.Range("J85").Value = 253
$ws.Range("J86").Value = 30
$ws.Range("J88").Value = 56
$ws.Range("J90").Value = 63
$ws.Range("J94").Value = 41
$ws.Range("J95").Value = 82
$ws.Range("J96").Value = 68
$ws.Range("J98").Value = 38
$ws.Range("J99").Value = 70
$ws.Range("I101").Value = 26201
$ws.Range("J101").Value = 5603

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 60
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 253

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 112
$ws.Range("J6").Value = 96
$ws.Range("J7").Value = 346

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 54
$ws.Range("J7").Value = 193

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 20
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 38
$ws.Range("J3").Value = 83
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 46
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 139

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 28
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 59
$ws.Range("J3").Value = 67
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 233

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 91
$ws.Range("J3").Value = 117
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 312

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 49
$ws.Range("J3").Value = 55
$ws.Range("J6").Value = 72
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 6
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 10
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J3").Value = 3
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 18
